$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column before EU (labeled "18-dec"),
#     shifting the "01-oct." .. "31-oct." columns one position to the right
#     (EU:FY -> EV:FZ), and fill the 24 data rows of the new column with "-".
$ws1 = $wb.Worksheets.Item("Prix Spot")
$ws1.Columns("EU:EU").Insert()
$ws1.Range("EU1").Value2 = "18-dec"
$ws1.Range("EU2:EU25").Value2 = "-"

# --- Sheet "Gaz": append a new row 181 with the 2025-12-16 price.
$ws2 = $wb.Worksheets.Item("Gaz")
$ws2.Range("A181").NumberFormat = "@"
$ws2.Range("A181").Value2 = "2025-12-16"
$ws2.Range("A181").Style = "Normal"
$ws2.Range("B181").Value2 = 25.475

# --- Sheet "CO2": append a new row 181 with the 2025-12-16 price.
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A181").NumberFormat = "@"
$ws3.Range("A181").Value2 = "2025-12-16"
$ws3.Range("A181").Style = "Normal"
$ws3.Range("B181").Value2 = 85.08
